# Auto-generated script applying scheduled market-data refresh to Sheets/Faerie_Profits.xlsx
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per leve row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 280.83334
$ws.Range("J12").Value = 811.5
$ws.Range("L12").Value = 811.5
$ws.Range("N12").Value = -1151.5
$ws.Range("H33").Value = 1814.1428
$ws.Range("I33").Value = 1739.8
$ws.Range("K33").Value = 1739.8
$ws.Range("M33").Value = -1510.8
$ws.Range("H52").Value = 6712.5
$ws.Range("J52").Value = 7000
$ws.Range("L52").Value = 21000
$ws.Range("N52").Value = -21320
$ws.Range("H96").Value = 300
$ws.Range("J96").Value = 200
$ws.Range("L96").Value = 600
$ws.Range("N96").Value = -3346
$ws.Range("H106").Value = 55705.5
$ws.Range("I106").Value = 76278.5
$ws.Range("J106").Value = 7701.8335
$ws.Range("K106").Value = 76278.5
$ws.Range("L106").Value = 7701.8335
$ws.Range("M106").Value = -75647.5
$ws.Range("N106").Value = -8963.833500000001
$ws.Range("H111").Value = 6653.9165
$ws.Range("I111").Value = 5984.8
$ws.Range("K111").Value = 17954.4
$ws.Range("M111").Value = -14887.4
$ws.Range("H132").Value = 3241.6667
$ws.Range("I132").Value = 3019.9092
$ws.Range("K132").Value = 9059.7276
$ws.Range("M132").Value = -6529.7276
$ws.Range("H137").Value = 3771.9312
$ws.Range("I137").Value = 4393.353
$ws.Range("K137").Value = 13180.059
$ws.Range("M137").Value = -10630.059
$ws.Range("H138").Value = 192882.1
$ws.Range("J138").Value = 217186.6
$ws.Range("L138").Value = 651559.8
$ws.Range("N138").Value = -661839.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3322.125
$ws.Range("I2").Value = 4345.25
$ws.Range("J2").Value = 2299
$ws.Range("K2").Value = 4345.25
$ws.Range("L2").Value = 2299
$ws.Range("M2").Value = -4232.25
$ws.Range("N2").Value = -2525
$ws.Range("H32").Value = 925.375
$ws.Range("I32").Value = 878.5402
$ws.Range("K32").Value = 878.5402
$ws.Range("M32").Value = -591.5402
$ws.Range("H45").Value = 4621.385
$ws.Range("I45").Value = 4154.2104
$ws.Range("J45").Value = 5889.4287
$ws.Range("K45").Value = 4154.2104
$ws.Range("L45").Value = 5889.4287
$ws.Range("M45").Value = -3777.2104
$ws.Range("N45").Value = -6643.4287
$ws.Range("H61").Value = 5093.8335
$ws.Range("I61").Value = 5201.763
$ws.Range("J61").Value = 4837.5
$ws.Range("K61").Value = 5201.763
$ws.Range("L61").Value = 4837.5
$ws.Range("M61").Value = -4989.763
$ws.Range("N61").Value = -5261.5
$ws.Range("H74").Value = 1729.2916
$ws.Range("I74").Value = 1884.2106
$ws.Range("J74").Value = 1140.6
$ws.Range("K74").Value = 1884.2106
$ws.Range("L74").Value = 1140.6
$ws.Range("M74").Value = -1010.2106
$ws.Range("N74").Value = -2888.6
$ws.Range("H77").Value = 1729.2916
$ws.Range("I77").Value = 1884.2106
$ws.Range("J77").Value = 1140.6
$ws.Range("K77").Value = 9421.053
$ws.Range("L77").Value = 5703
$ws.Range("M77").Value = -5053.053
$ws.Range("N77").Value = -14439
$ws.Range("H88").Value = 3282
$ws.Range("I88").Value = 3500
$ws.Range("J88").Value = 3194.8
$ws.Range("K88").Value = 3500
$ws.Range("L88").Value = 3194.8
$ws.Range("M88").Value = -3094
$ws.Range("N88").Value = -4006.8
$ws.Range("H91").Value = 3282
$ws.Range("I91").Value = 3500
$ws.Range("J91").Value = 3194.8
$ws.Range("K91").Value = 3500
$ws.Range("L91").Value = 3194.8
$ws.Range("M91").Value = -2096
$ws.Range("N91").Value = -6002.8
$ws.Range("H98").Value = 94999.336
$ws.Range("J98").Value = 94999.336
$ws.Range("L98").Value = 94999.336
$ws.Range("N98").Value = -100989.336
$ws.Range("H116").Value = 3322.125
$ws.Range("I116").Value = 4345.25
$ws.Range("J116").Value = 2299
$ws.Range("K116").Value = 4345.25
$ws.Range("L116").Value = 2299
$ws.Range("M116").Value = -2051.25
$ws.Range("N116").Value = -6887
$ws.Range("H122").Value = 1596.3055
$ws.Range("I122").Value = 1353.1666
$ws.Range("K122").Value = 4059.4998
$ws.Range("M122").Value = -1609.4998
$ws.Range("H132").Value = 3842.2563
$ws.Range("I132").Value = 2261.9048
$ws.Range("J132").Value = 5686
$ws.Range("K132").Value = 6785.714399999999
$ws.Range("L132").Value = 17058
$ws.Range("M132").Value = -4255.714399999999
$ws.Range("N132").Value = -22118
$ws.Range("H136").Value = 5093.8335
$ws.Range("I136").Value = 5201.763
$ws.Range("J136").Value = 4837.5
$ws.Range("K136").Value = 15605.289
$ws.Range("L136").Value = 14512.5
$ws.Range("M136").Value = -13055.289
$ws.Range("N136").Value = -19612.5
$ws.Range("H141").Value = 67214.5
$ws.Range("J141").Value = 67214.5
$ws.Range("L141").Value = 67214.5
$ws.Range("N141").Value = -77574.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3322.125
$ws.Range("I3").Value = 4345.25
$ws.Range("J3").Value = 2299
$ws.Range("K3").Value = 4345.25
$ws.Range("L3").Value = 2299
$ws.Range("M3").Value = -4231.25
$ws.Range("N3").Value = -2527
$ws.Range("H20").Value = 985.04
$ws.Range("I20").Value = 811.64703
$ws.Range("K20").Value = 811.64703
$ws.Range("M20").Value = -564.64703
$ws.Range("H22").Value = 433.33334
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -127
$ws.Range("N22").Value = -846
$ws.Range("H80").Value = 393.33334
$ws.Range("I80").Value = 197.14285
$ws.Range("J80").Value = 491.42856
$ws.Range("K80").Value = 197.14285
$ws.Range("L80").Value = 491.42856
$ws.Range("M80").Value = 800.85715
$ws.Range("N80").Value = -2487.42856
$ws.Range("H83").Value = 393.33334
$ws.Range("I83").Value = 197.14285
$ws.Range("J83").Value = 491.42856
$ws.Range("K83").Value = 985.71425
$ws.Range("L83").Value = 2457.1428
$ws.Range("M83").Value = 4006.28575
$ws.Range("N83").Value = -12441.1428
$ws.Range("H94").Value = 3327.4614
$ws.Range("I94").Value = 2889.7144
$ws.Range("J94").Value = 3838.1667
$ws.Range("K94").Value = 2889.7144
$ws.Range("L94").Value = 3838.1667
$ws.Range("M94").Value = -2438.7144
$ws.Range("N94").Value = -4740.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 8388333.5
$ws.Range("I6").Value = 12506250
$ws.Range("J6").Value = 152500
$ws.Range("K6").Value = 12506250
$ws.Range("L6").Value = 152500
$ws.Range("M6").Value = -12506137
$ws.Range("N6").Value = -152726
$ws.Range("H11").Value = 290
$ws.Range("I11").Value = 290
$ws.Range("K11").Value = 290
$ws.Range("M11").Value = -150
$ws.Range("H22").Value = 541.125
$ws.Range("I22").Value = 553.5833
$ws.Range("K22").Value = 553.5833
$ws.Range("M22").Value = -203.5833
$ws.Range("H31").Value = 3343.5938
$ws.Range("I31").Value = 2441.2354
$ws.Range("J31").Value = 4366.2666
$ws.Range("K31").Value = 2441.2354
$ws.Range("L31").Value = 4366.2666
$ws.Range("M31").Value = -2146.2354
$ws.Range("N31").Value = -4956.2666
$ws.Range("H32").Value = 8957.4
$ws.Range("I32").Value = 8957.4
$ws.Range("K32").Value = 8957.4
$ws.Range("M32").Value = -8641.4
$ws.Range("H34").Value = 3343.5938
$ws.Range("I34").Value = 2441.2354
$ws.Range("J34").Value = 4366.2666
$ws.Range("K34").Value = 2441.2354
$ws.Range("L34").Value = 4366.2666
$ws.Range("M34").Value = -2239.2354
$ws.Range("N34").Value = -4770.2666
$ws.Range("H92").Value = 71966.664
$ws.Range("J92").Value = 71966.664
$ws.Range("L92").Value = 71966.664
$ws.Range("N92").Value = -76958.664
$ws.Range("H122").Value = 1270.909
$ws.Range("I122").Value = 970
$ws.Range("J122").Value = 2625
$ws.Range("K122").Value = 2910
$ws.Range("L122").Value = 7875
$ws.Range("M122").Value = -460
$ws.Range("N122").Value = -12775
$ws.Range("H132").Value = 2498.25
$ws.Range("I132").Value = 2498.25
$ws.Range("K132").Value = 7494.75
$ws.Range("M132").Value = -4964.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 456.8
$ws.Range("I12").Value = 979.8
$ws.Range("J12").Value = 282.46667
$ws.Range("K12").Value = 2939.4
$ws.Range("L12").Value = 847.4000100000001
$ws.Range("M12").Value = -2766.4
$ws.Range("N12").Value = -1193.40001
$ws.Range("H59").Value = 97.5
$ws.Range("I59").Value = 97.5
$ws.Range("K59").Value = 292.5
$ws.Range("M59").Value = 247.5
$ws.Range("H68").Value = 1797.8422
$ws.Range("J68").Value = 1913.5
$ws.Range("L68").Value = 5740.5
$ws.Range("N68").Value = -7362.5
$ws.Range("H71").Value = 1797.8422
$ws.Range("J71").Value = 1913.5
$ws.Range("L71").Value = 17221.5
$ws.Range("N71").Value = -25333.5
$ws.Range("H80").Value = 5094.8
$ws.Range("I80").Value = 2247
$ws.Range("J80").Value = 6993.3335
$ws.Range("K80").Value = 6741
$ws.Range("L80").Value = 20980.0005
$ws.Range("M80").Value = -5805
$ws.Range("N80").Value = -22852.0005
$ws.Range("H83").Value = 5094.8
$ws.Range("I83").Value = 2247
$ws.Range("J83").Value = 6993.3335
$ws.Range("K83").Value = 20223
$ws.Range("L83").Value = 62940.0015
$ws.Range("M83").Value = -15543
$ws.Range("N83").Value = -72300.0015
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 6000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -7622
$ws.Range("H113").Value = 586.4666999999999
$ws.Range("I113").Value = 354.14285
$ws.Range("J113").Value = 789.75
$ws.Range("K113").Value = 1062.42855
$ws.Range("L113").Value = 2369.25
$ws.Range("M113").Value = 1107.57145
$ws.Range("N113").Value = -6709.25
$ws.Range("H122").Value = 924.1
$ws.Range("I122").Value = 678.6
$ws.Range("J122").Value = 1169.6
$ws.Range("K122").Value = 6107.400000000001
$ws.Range("L122").Value = 10526.4
$ws.Range("M122").Value = -3657.400000000001
$ws.Range("N122").Value = -15426.4
$ws.Range("H125").Value = 8599.666999999999
$ws.Range("I125").Value = 5799
$ws.Range("J125").Value = 10000
$ws.Range("K125").Value = 17397
$ws.Range("L125").Value = 30000
$ws.Range("M125").Value = -12477
$ws.Range("N125").Value = -39840
$ws.Range("H127").Value = 1727.7
$ws.Range("J127").Value = 1727.7
$ws.Range("L127").Value = 5183.1
$ws.Range("N127").Value = -15103.1
$ws.Range("H129").Value = 1788.625
$ws.Range("I129").Value = 444.8
$ws.Range("J129").Value = 8507.75
$ws.Range("K129").Value = 1334.4
$ws.Range("L129").Value = 25523.25
$ws.Range("M129").Value = 3665.6
$ws.Range("N129").Value = -35523.25
$ws.Range("H130").Value = 3998.5
$ws.Range("I130").Value = 3998.5
$ws.Range("K130").Value = 11995.5
$ws.Range("M130").Value = -6975.5
$ws.Range("H131").Value = 3823.484
$ws.Range("I131").Value = 6311.7856
$ws.Range("J131").Value = 1774.2941
$ws.Range("K131").Value = 18935.3568
$ws.Range("L131").Value = 5322.8823
$ws.Range("M131").Value = -13895.3568
$ws.Range("N131").Value = -15402.8823
$ws.Range("H133").Value = 6998.4287
$ws.Range("I133").Value = 2997.3333
$ws.Range("J133").Value = 9999.25
$ws.Range("K133").Value = 8991.999899999999
$ws.Range("L133").Value = 29997.75
$ws.Range("M133").Value = -3931.999899999999
$ws.Range("N133").Value = -40117.75
$ws.Range("H137").Value = 13893770
$ws.Range("J137").Value = 25647930
$ws.Range("L137").Value = 76943790
$ws.Range("N137").Value = -76953990

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5603.316
$ws.Range("I70").Value = 2962.7856
$ws.Range("K70").Value = 2962.7856
$ws.Range("M70").Value = -2692.7856
$ws.Range("H73").Value = 5603.316
$ws.Range("I73").Value = 2962.7856
$ws.Range("K73").Value = 2962.7856
$ws.Range("M73").Value = -2026.7856
$ws.Range("H80").Value = 12759
$ws.Range("I80").Value = 1275
$ws.Range("J80").Value = 20415
$ws.Range("K80").Value = 1275
$ws.Range("L80").Value = 20415
$ws.Range("M80").Value = -277
$ws.Range("N80").Value = -22411
$ws.Range("H83").Value = 12759
$ws.Range("I83").Value = 1275
$ws.Range("J83").Value = 20415
$ws.Range("K83").Value = 6375
$ws.Range("L83").Value = 102075
$ws.Range("M83").Value = -1383
$ws.Range("N83").Value = -112059
$ws.Range("H92").Value = 20386.2
$ws.Range("J92").Value = 20386.2
$ws.Range("L92").Value = 20386.2
$ws.Range("N92").Value = -24130.2
$ws.Range("H93").Value = 83532
$ws.Range("J93").Value = 83532
$ws.Range("L93").Value = 83532
$ws.Range("N93").Value = -87276
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 2184.9167
$ws.Range("I97").Value = 1141.9
$ws.Range("J97").Value = 7400
$ws.Range("K97").Value = 1141.9
$ws.Range("L97").Value = 7400
$ws.Range("M97").Value = -645.9000000000001
$ws.Range("N97").Value = -8392
$ws.Range("H102").Value = 24709.455
$ws.Range("I102").Value = 1299.7646
$ws.Range("K102").Value = 1299.7646
$ws.Range("M102").Value = 322.2354
$ws.Range("H124").Value = 87000
$ws.Range("J124").Value = 87000
$ws.Range("L124").Value = 87000
$ws.Range("N124").Value = -96820
$ws.Range("H126").Value = 6626.353
$ws.Range("I126").Value = 3382.5
$ws.Range("K126").Value = 10147.5
$ws.Range("M126").Value = -7677.5
$ws.Range("H132").Value = 5408.36
$ws.Range("I132").Value = 6165.15
$ws.Range("K132").Value = 18495.45
$ws.Range("M132").Value = -15965.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 20009
$ws.Range("I5").Value = 20009
$ws.Range("K5").Value = 20009
$ws.Range("M5").Value = -19896
$ws.Range("H22").Value = 3633.1333
$ws.Range("I22").Value = 2360.8
$ws.Range("K22").Value = 2360.8
$ws.Range("M22").Value = -2065.8
$ws.Range("H27").Value = 3633.1333
$ws.Range("I27").Value = 2360.8
$ws.Range("K27").Value = 2360.8
$ws.Range("M27").Value = -2253.8
$ws.Range("H40").Value = 1798.7142
$ws.Range("I40").Value = 1798.7142
$ws.Range("K40").Value = 1798.7142
$ws.Range("M40").Value = -1662.7142
$ws.Range("H93").Value = 1376.6
$ws.Range("I93").Value = 1376.6
$ws.Range("K93").Value = 1376.6
$ws.Range("M93").Value = -128.5999999999999
$ws.Range("H100").Value = 4444.3335
$ws.Range("I100").Value = 3999.8
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 3999.8
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -3458.8
$ws.Range("N100").Value = -6082
$ws.Range("H122").Value = 5262.8203
$ws.Range("I122").Value = 3795.25
$ws.Range("K122").Value = 11385.75
$ws.Range("M122").Value = -8935.75
$ws.Range("H132").Value = 7701.4634
$ws.Range("I132").Value = 7791.593
$ws.Range("J132").Value = 7527.643
$ws.Range("K132").Value = 23374.779
$ws.Range("L132").Value = 22582.929
$ws.Range("M132").Value = -20844.779
$ws.Range("N132").Value = -27642.929
$ws.Range("H138").Value = 78214.5
$ws.Range("J138").Value = 78214.5
$ws.Range("L138").Value = 78214.5
$ws.Range("N138").Value = -88494.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 19978.75
$ws.Range("I21").Value = 21671.666
$ws.Range("J21").Value = 14900
$ws.Range("K21").Value = 21671.666
$ws.Range("L21").Value = 14900
$ws.Range("M21").Value = -21436.666
$ws.Range("N21").Value = -15370
$ws.Range("H28").Value = 26086.166
$ws.Range("J28").Value = 28375
$ws.Range("L28").Value = 28375
$ws.Range("N28").Value = -29071
$ws.Range("H30").Value = 10331.777
$ws.Range("J30").Value = 8329.5
$ws.Range("L30").Value = 8329.5
$ws.Range("N30").Value = -8543.5
$ws.Range("H35").Value = 19978.75
$ws.Range("I35").Value = 21671.666
$ws.Range("J35").Value = 14900
$ws.Range("K35").Value = 21671.666
$ws.Range("L35").Value = 14900
$ws.Range("M35").Value = -21381.666
$ws.Range("N35").Value = -15480
$ws.Range("H100").Value = 1279.7693
$ws.Range("J100").Value = 1448.875
$ws.Range("L100").Value = 2897.75
$ws.Range("N100").Value = -3979.75
$ws.Range("H107").Value = 567.3461
$ws.Range("I107").Value = 489.26666
$ws.Range("J107").Value = 673.8182
$ws.Range("K107").Value = 1467.79998
$ws.Range("L107").Value = 2021.4546
$ws.Range("M107").Value = 452.20002
$ws.Range("N107").Value = -5861.4546
$ws.Range("H113").Value = 4166969.5
$ws.Range("I113").Value = 5555911.5
$ws.Range("J113").Value = 143.2
$ws.Range("K113").Value = 16667734.5
$ws.Range("L113").Value = 429.6
$ws.Range("M113").Value = -16665564.5
$ws.Range("N113").Value = -4769.6
$ws.Range("H122").Value = 7371.1763
$ws.Range("I122").Value = 5370.4287
$ws.Range("J122").Value = 8771.700000000001
$ws.Range("K122").Value = 16111.2861
$ws.Range("L122").Value = 26315.1
$ws.Range("M122").Value = -13661.2861
$ws.Range("N122").Value = -31215.1
$ws.Range("H126").Value = 2183.439
$ws.Range("I126").Value = 2164.9062
$ws.Range("K126").Value = 6494.7186
$ws.Range("M126").Value = -4024.7186
$ws.Range("H132").Value = 2610.6765
$ws.Range("I132").Value = 1902.5
$ws.Range("K132").Value = 5707.5
$ws.Range("M132").Value = -3177.5
$ws.Range("H136").Value = 2908.5312
$ws.Range("I136").Value = 1438.3684
$ws.Range("J136").Value = 5057.231
$ws.Range("K136").Value = 4315.1052
$ws.Range("L136").Value = 15171.693
$ws.Range("M136").Value = -1765.1052
$ws.Range("N136").Value = -20271.693

